# Highlight best experimental results in the excel files
$wb = $excel.ActiveWorkbook

# Rows/columns for the three result blocks on each sheet.
# For each block: row of "r" (correlation -> higher better, MAX)
# and rows of "MAE"/"RMSE" (lower better, MIN)
$maxRows = @(2, 8, 14)
$minRows = @(3, 4, 9, 10, 15, 16)

foreach ($sheetName in @("Failed", "Survived")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($r in $maxRows) {
        $src = "B" + $r + ":E" + $r
        $dst = $ws.Range("G" + $r)
        $dst.Formula = "=MAX(" + $src + ")"

        $rng = $ws.Range("B" + $r + ":E" + $r)
        $best = $rng.Cells.Item(1, 1)
        for ($c = 1; $c -le 4; $c++) {
            $cell = $rng.Cells.Item(1, $c)
            if ($cell.Value -gt $best.Value) { $best = $cell }
        }
        $best.Font.Bold = $true
        $best.Font.Color = 255
    }

    foreach ($r in $minRows) {
        $src = "B" + $r + ":E" + $r
        $dst = $ws.Range("G" + $r)
        $dst.Formula = "=MIN(" + $src + ")"

        $rng = $ws.Range("B" + $r + ":E" + $r)
        $best = $rng.Cells.Item(1, 1)
        for ($c = 1; $c -le 4; $c++) {
            $cell = $rng.Cells.Item(1, $c)
            if ($cell.Value -lt $best.Value) { $best = $cell }
        }
        $best.Font.Bold = $true
        $best.Font.Color = 255
    }
}

# Page setup on "Failed" sheet (A4, portrait)
$wsFailed = $wb.Worksheets.Item("Failed")
$wsFailed.PageSetup.PaperSize = 9
$wsFailed.PageSetup.Orientation = 1

# Selections
$wsFailed.Range("C21").Select()
$wsSurvived = $wb.Worksheets.Item("Survived")
$wsSurvived.Range("M8").Select()

# Make "Survived" the active/tab-selected sheet
$wsSurvived.Activate()
